# RefApp-Framework - added information about header files, minor corrections.
#
# 1) Bump the cached "today" date shown in the Handout Master and the
#    Notes Master footers from 20/02/2023 to 21/02/2023.
# 2) Slide 14 (the header-pin-out diagram): fix the pin labels.
#    - The first I2C-labelled connector was mislabelled UART pins;
#      rename its SDA/SCL boxes to RXD/TXD.
#    - The third connector (under the SPI label) had its SDA/SCL
#      labels swapped; correct them.

$p = $ppt.ActivePresentation

# --- Handout master / Notes master date placeholders ---------------------
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = "21/02/2023"
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "21/02/2023"

# --- Slide 14 pin labels ---------------------------------------------------
$s = $p.Slides.Item(14)

$s.Shapes.Item("Rectangle 102").TextFrame.TextRange.Text = "RXD"
$s.Shapes.Item("Rectangle 103").TextFrame.TextRange.Text = "TXD"

$s.Shapes.Item("Rectangle 155").TextFrame.TextRange.Text = "SCL"
$s.Shapes.Item("Rectangle 156").TextFrame.TextRange.Text = "SDA"
